$d = $word.ActiveDocument

$lsquo = [char]0x2018
$rsquo = [char]0x2019

# 1) "There do not ... differences between Transrate scores in the NCGR 'nt' asse"
#    -> "Frequency distribution of differences between Transrate scores between the NCGR 'nt' asse"
#    (one single Find/Replace so the match fully spans both proofErr-wrapped runs -
#     "Transrate" and "nt" - and those stray proofErr markers get dropped)
$d.Content.Find.Execute(
    "There do not appear to be taxonomic trends in the frequency of differences between Transrate scores in the NCGR " + $lsquo + "nt" + $rsquo + " asse",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Frequency distribution of differences between Transrate scores between the NCGR " + $lsquo + "nt" + $rsquo + " asse", 2)

# 3) "mblies and the DIB re-assemblies." -> "mblies and the DIB re-assemblies grouped by the top seven most represented phyla in the MMETSP data set."
$d.Content.Find.Execute(
    "mblies and the DIB re-assemblies.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "mblies and the DIB re-assemblies grouped by the top seven most represented phyla in the MMETSP data set.", 2)

# 4) Rewrite the "Negative values indicate ... were higher." sentence, in two pieces that
#    stop/resume exactly at the "re-assemblies" / " were higher." boundary so the
#    "_GoBack" bookmark sitting there is left alone (a Find/Replace spanning across a
#    bookmark would delete it).
#    (the first piece fully spans the Transrate proofErr-wrapped run so its stray
#     proofErr markers are dropped)
$d.Content.Find.Execute(
    "Negative values indicate Transrate scores from NCGR assemblies were higher and positive values indicate scores from DIB re-assemblies",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Negative values indicate that Transrate scores from NCGR " + $lsquo + "nt" + $rsquo + " assemblies were higher than the score from the DIB re-assemblies and positive values indicate that scores from DIB re-assemblies", 2)

$d.Content.Find.Execute(
    " were higher.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " were higher than the NCGR " + $lsquo + "nt" + $rsquo + " assemblies.", 2)
